$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update the greeting text in E8 from "Good Morning" to "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Reflect the active selection on the sheet being cell E8
$ws.Range("E8").Select()
